$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the two invoice rows for 64/23-24 and 66/23-24 (rows 2 and 3).
# This shifts all subsequent rows up by two.
$ws.Rows.Item(2).Delete()
$ws.Rows.Item(2).Delete()

# Restore the "Sr. No" label that was on the deleted row 2 (64/23-24).
$ws.Range("A2").Value = 1

# Update the sum formula for the remaining two rows in the first group.
$ws.Range("F3").Formula = "=E2+E3"

# Add the small annotation row (single space) just below the first group.
$ws.Range("F4").Value = " "
